# Auto-generated edit script: update Leve profit calculation sheets
# with refreshed market data (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 463654.94
$ws.Range("J17").Value = 514998.2
$ws.Range("L17").Value = 1544994.6
$ws.Range("N17").Value = -1545330.6
$ws.Range("H33").Value = 4347977.5
$ws.Range("I33").Value = 5263289
$ws.Range("K33").Value = 5263289
$ws.Range("M33").Value = -5263060
$ws.Range("H40").Value = 1203.6897
$ws.Range("J40").Value = 1231.3077
$ws.Range("L40").Value = 1231.3077
$ws.Range("N40").Value = -1581.3077
$ws.Range("H96").Value = 908
$ws.Range("I96").Value = 884
$ws.Range("K96").Value = 2652
$ws.Range("M96").Value = -1279
$ws.Range("H100").Value = 2072
$ws.Range("I100").Value = 2143.5
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 2143.5
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -1602.5
$ws.Range("N100").Value = -2582
$ws.Range("H101").Value = 328.42856
$ws.Range("I101").Value = 341.5
$ws.Range("K101").Value = 1024.5
$ws.Range("M101").Value = 597.5
$ws.Range("H137").Value = 1028.6285
$ws.Range("I137").Value = 994.4516
$ws.Range("K137").Value = 2983.3548
$ws.Range("M137").Value = -433.3548000000001
$ws.Range("H141").Value = 1771.25
$ws.Range("I141").Value = 1310
$ws.Range("K141").Value = 3930
$ws.Range("M141").Value = 1250

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 36335.332
$ws.Range("J23").Value = 36335.332
$ws.Range("L23").Value = 36335.332
$ws.Range("N23").Value = -36853.332
$ws.Range("H32").Value = 2270.3293
$ws.Range("I32").Value = 1849.581
$ws.Range("K32").Value = 1849.581
$ws.Range("M32").Value = -1562.581
$ws.Range("H38").Value = 18491.727
$ws.Range("I38").Value = 5376.4443
$ws.Range("J38").Value = 77510.5
$ws.Range("K38").Value = 5376.4443
$ws.Range("L38").Value = 77510.5
$ws.Range("M38").Value = -4909.4443
$ws.Range("N38").Value = -78444.5
$ws.Range("H74").Value = 5554.476
$ws.Range("I74").Value = 619.1667
$ws.Range("K74").Value = 619.1667
$ws.Range("M74").Value = 254.8333
$ws.Range("H77").Value = 5554.476
$ws.Range("I77").Value = 619.1667
$ws.Range("K77").Value = 3095.8335
$ws.Range("M77").Value = 1272.1665
$ws.Range("H102").Value = 3772
$ws.Range("I102").Value = 3650.6667
$ws.Range("K102").Value = 3650.6667
$ws.Range("M102").Value = -2028.6667
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 2922
$ws.Range("I122").Value = 2444
$ws.Range("K122").Value = 7332
$ws.Range("M122").Value = -4882
$ws.Range("H132").Value = 2089.4614
$ws.Range("I132").Value = 1893.04
$ws.Range("K132").Value = 5679.12
$ws.Range("M132").Value = -3149.12

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2030.0952
$ws.Range("I20").Value = 1721.5518
$ws.Range("K20").Value = 1721.5518
$ws.Range("M20").Value = -1474.5518
$ws.Range("H30").Value = 15011
$ws.Range("J30").Value = 15011
$ws.Range("L30").Value = 15011
$ws.Range("N30").Value = -15261
$ws.Range("H94").Value = 1352.75
$ws.Range("J94").Value = 3250
$ws.Range("L94").Value = 3250
$ws.Range("N94").Value = -4152

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1288.6471
$ws.Range("I94").Value = 1390.6666
$ws.Range("K94").Value = 1390.6666
$ws.Range("M94").Value = -939.6666
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H131").Value = 27542.188
$ws.Range("I131").Value = 30000
$ws.Range("J131").Value = 27191.072
$ws.Range("K131").Value = 30000
$ws.Range("L131").Value = 27191.072
$ws.Range("M131").Value = -24960
$ws.Range("N131").Value = -37271.072

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 124.35714
$ws.Range("J12").Value = 141.6
$ws.Range("L12").Value = 424.8
$ws.Range("N12").Value = -770.8
$ws.Range("H39").Value = 4891.5386
$ws.Range("J39").Value = 4891.5386
$ws.Range("L39").Value = 14674.6158
$ws.Range("N39").Value = -15262.6158
$ws.Range("H56").Value = 13994.6
$ws.Range("I56").Value = 13994.6
$ws.Range("K56").Value = 13994.6
$ws.Range("M56").Value = -13464.6
$ws.Range("H92").Value = 587.0476
$ws.Range("I92").Value = 269.27274
$ws.Range("J92").Value = 936.6
$ws.Range("K92").Value = 807.81822
$ws.Range("L92").Value = 2809.8
$ws.Range("M92").Value = 440.18178
$ws.Range("N92").Value = -5305.8
$ws.Range("H134").Value = 5640.2
$ws.Range("I134").Value = 3133.5833
$ws.Range("J134").Value = 15666.667
$ws.Range("K134").Value = 9400.749899999999
$ws.Range("L134").Value = 47000.001
$ws.Range("M134").Value = -4330.749899999999
$ws.Range("N134").Value = -57140.001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1591.6522
$ws.Range("J102").Value = 4325
$ws.Range("L102").Value = 4325
$ws.Range("N102").Value = -7569
$ws.Range("H122").Value = 2546.5386
$ws.Range("I122").Value = 2410.4546
$ws.Range("K122").Value = 7231.3638
$ws.Range("M122").Value = -4781.3638

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3552.311
$ws.Range("I16").Value = 1236.04
$ws.Range("K16").Value = 1236.04
$ws.Range("M16").Value = -1066.04
$ws.Range("H55").Value = 161.42105
$ws.Range("I55").Value = 160.5
$ws.Range("J55").Value = 166.33333
$ws.Range("K55").Value = 160.5
$ws.Range("L55").Value = 166.33333
$ws.Range("M55").Value = 12.5
$ws.Range("N55").Value = -512.3333299999999
$ws.Range("H93").Value = 587668.1
$ws.Range("I93").Value = 857631.5600000001
$ws.Range("K93").Value = 857631.5600000001
$ws.Range("M93").Value = -856383.5600000001
$ws.Range("H100").Value = 73982.47
$ws.Range("I100").Value = 76838.36
$ws.Range("K100").Value = 76838.36
$ws.Range("M100").Value = -76297.36
$ws.Range("H122").Value = 10572.728
$ws.Range("I122").Value = 9400.166999999999
$ws.Range("K122").Value = 28200.501
$ws.Range("M122").Value = -25750.501
$ws.Range("H132").Value = 4059.1875
$ws.Range("I132").Value = 2919
$ws.Range("K132").Value = 8757
$ws.Range("M132").Value = -6227

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 29002
$ws.Range("I9").Value = 32499.5
$ws.Range("J9").Value = 22007
$ws.Range("K9").Value = 32499.5
$ws.Range("L9").Value = 22007
$ws.Range("M9").Value = -32359.5
$ws.Range("N9").Value = -22287
$ws.Range("H100").Value = 979.4595
$ws.Range("J100").Value = 1276.8182
$ws.Range("L100").Value = 2553.6364
$ws.Range("N100").Value = -3635.6364
$ws.Range("H122").Value = 4132.0625
$ws.Range("I122").Value = 3959.5
$ws.Range("J122").Value = 4304.625
$ws.Range("K122").Value = 11878.5
$ws.Range("L122").Value = 12913.875
$ws.Range("M122").Value = -9428.5
$ws.Range("N122").Value = -17813.875
$ws.Range("H132").Value = 2498.9424
$ws.Range("I132").Value = 2808.442
$ws.Range("K132").Value = 8425.326000000001
$ws.Range("M132").Value = -5895.326000000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

